# Commit: "Rajout de la config JPA et connexion a la base"
#
# The corresponding OOXML diff removes the 4th slide of the deck
# (sldId="260", r:id="rId5" -> ppt/slides/slide4.xml, title
# "https://start.spring.io/") from the slide list. PowerPoint's Slides
# collection is 1-based and in this deck slide index 4 is exactly that
# slide, so deleting it removes the <p:sldId .../> entry from
# <p:sldIdLst>, drops ppt/slides/slide4.xml + its _rels/image, and
# removes the matching Content_Types override / presentation.xml.rels
# relationship automatically.

$p = $ppt.ActivePresentation

$target = $null
foreach ($s in $p.Slides) {
    if ($s.SlideID -eq 260) {
        $target = $s
        break
    }
}

if ($target -eq $null) {
    # Fallback: locate by the slide's title text if the SlideID ever
    # differs from the one recorded in the diff.
    foreach ($s in $p.Slides) {
        foreach ($sh in $s.Shapes) {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "https://start.spring.io/") {
                    $target = $s
                }
                break
            }
        }
        if ($target -ne $null) { break }
    }
}

if ($target -eq $null) {
    # Last resort: the slide to remove is the last one in this deck.
    $target = $p.Slides.Item($p.Slides.Count)
}

$target.Delete()
